$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.478.23"
$ws.Range("E2").Value = "  +1.73%  "

$ws.Range("D3").Value = "1.880.03"
$ws.Range("E3").Value = "  +2.01%  "

$ws.Range("E4").Value = "  +0.41%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.88%  "

$ws.Range("E7").Value = "  +0.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.333"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0699"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0991"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.06%  "

$ws.Range("D12").Value = "2.145.28"
$ws.Range("E12").Value = "  +1.72%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.687"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.78%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.23%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.822.42"
$ws.Range("E16").Value = "  -0.89%  "

$ws.Range("D17").Value = "35.512.87"
$ws.Range("E17").Value = "  +1.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.82%  "

$ws.Range("D19").Value = "0.0₃0803"
$ws.Range("E19").Value = "  +1.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "242.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.62%  "

$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +26.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.07%  "

$ws.Range("E29").Value = "  +1.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0565"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.23%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.54%  "

$ws.Range("B32").Value = "BinanceUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +24.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.830"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +18.96%  "

$ws.Range("E37").Value = "  +6.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.54%  "

$ws.Range("E39").Value = "  +4.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "91.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.24%  "

$ws.Range("D41").Value = "1.355.89"
$ws.Range("E41").Value = "  +0.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.94%  "

$ws.Range("E43").Value = "  +14.92%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.30%  "

$ws.Range("B45").Value = "Gas"
$ws.Range("C45").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +58.14%  "

$ws.Range("E46").Value = "  +0.31%  "

$ws.Range("E47").Value = "  +6.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.64%  "

$ws.Range("D49").Value = "2.063.34"
$ws.Range("E49").Value = "  +2.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0689"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.54%  "
